$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169; this pushes the existing rows 169-269 down to 170-270.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new record's data.
$ws.Cells.Item(169, 1).Value  = 4
$ws.Cells.Item(169, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(169, 3).Value  = "Los Lagos"
$ws.Cells.Item(169, 4).Value  = 44596
$ws.Cells.Item(169, 5).Value  = 10
$ws.Cells.Item(169, 6).Value  = 100112008
$ws.Cells.Item(169, 7).Value  = "Coliflor"
$ws.Cells.Item(169, 8).Value  = "Sin especificar"
$ws.Cells.Item(169, 9).Value  = "Primera"
$ws.Cells.Item(169, 10).Value = 600
$ws.Cells.Item(169, 11).Value = 1500
$ws.Cells.Item(169, 12).Value = 1600
$ws.Cells.Item(169, 13).Value = 1550
$ws.Cells.Item(169, 14).Value = "$/unidad"
$ws.Cells.Item(169, 15).Value = "Región Metropolitana"
$ws.Cells.Item(169, 16).Value = 1550
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"
